$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-10 (row 11 will be removed afterwards)
$data = @(
    @("1331470", "https://aiesec.org/opportunity/global-talent/1331470", "CONSTRUCTION PROJECT COORDINATOR", "Abidjan, Côte d'Ivoire", "No", "0 applicants", "9 - 12 Weeks", "ONG REFUGE DES ENFANTS"),
    @("1331382", "https://aiesec.org/opportunity/global-talent/1331382", "Sales & Marketing Intern", "Pune, Maharashtra, India", "No", "1 applicant", "6 - 18 Months", "All Home Living LLP"),
    @("1331380", "https://aiesec.org/opportunity/global-talent/1331380", "Business Operations & Community Development Intern", "Jaipur, Rajasthan, India", "No", "0 applicants", "3 - 6 Months", "Boomerang Hospitality"),
    @("1330067", "https://aiesec.org/opportunity/global-talent/1330067", "Associate Program Manager (Eu Citizen Only)", "Bruxelles, Belgio", "No", "38 applicants", "6 - 18 Months", "The Synergist"),
    @("1329617", "https://aiesec.org/opportunity/global-talent/1329617", "[Impact Brazil] - Comercial and Sales Consultant Intern", "São Paulo, SP, Brasil", "No", "59 applicants", "6 - 18 Months", "Hays Brasil"),
    @("1328614", "https://aiesec.org/opportunity/global-talent/1328614", "Field Service Engineer [EU Preferred]", "Madrid, Spain", "No", "159 applicants", "6 - 18 Months", "Mitsubishi Power Europe Sucursal en España"),
    @("1327183", "https://aiesec.org/opportunity/global-talent/1327183", "Sales ENGINEERING", "Çayırova, Kocaeli, Türkiye", "No", "113 applicants", "3 - 6 Months", "Eurotray Metal Elektrik San. Tic.: Ltd. Şti."),
    @("1316099", "https://aiesec.org/opportunity/global-talent/1316099", "Sales & BD Junior Manager", "Puzi City, Chiayi County, Taiwan 613", "No", "100 applicants", "6 - 18 Months", "Singform Enterprise Co., Ltd."),
    @("1296457", "https://aiesec.org/opportunity/global-talent/1296457", "Export Sales Specialist in EUROTRAY Company", "Dilovası, Türkiye", "No", "179 applicants", "6 - 18 Months", "Eurotray Metal Elektrik Company")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]

    # Column A holds numeric-looking opportunity IDs; force text so Excel
    # doesn't silently convert them to numbers (matches source inlineStr text).
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 1).Style = "Normal"

    for ($c = 2; $c -le 8; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}

# Remove the now-obsolete row 11 (previously the 10th data row), shifting cells up
$ws.Rows.Item(11).Delete()

# Adjust column widths per the updated layout.
# The COM ColumnWidth setter adds a constant ~5/6 character padding versus
# the raw OOXML "width" attribute, so compensate to land on the exact
# target widths (58, 39, 17, 47).
$pad = 0.8333333333333334
$ws.Columns.Item(3).ColumnWidth = 58 - $pad
$ws.Columns.Item(4).ColumnWidth = 39 - $pad
$ws.Columns.Item(6).ColumnWidth = 17 - $pad
$ws.Columns.Item(8).ColumnWidth = 47 - $pad
